{"js": "// Apply the text replacements described by the diff: update the date line\n// and each two-digit multiplication answer cell in the table.\nconst replacements = [\n  [\"2023-08-21 Monday\", \"2023-08-22 Tuesday\"],\n  [\"41\u00d722=902\", \"75\u00d759=4425\"],\n  [\"14\u00d750=700\", \"97\u00d732=3104\"],\n  [\"61\u00d766=4026\", \"75\u00d787=6525\"],\n  [\"44\u00d726=1144\", \"15\u00d747=705\"],\n  [\"41\u00d725=1025\", \"20\u00d733=660\"],\n  [\"87\u00d752=4524\", \"37\u00d762=2294\"],\n  [\"15\u00d757=855\", \"64\u00d780=5120\"],\n  [\"28\u00d713=364\", \"53\u00d799=5247\"],\n  [\"36\u00d728=1008\", \"29\u00d720=580\"],\n  [\"94\u00d791=8554\", \"95\u00d796=9120\"],\n  [\"21\u00d789=1869\", \"14\u00d749=686\"],\n  [\"64\u00d785=5440\", \"92\u00d749=4508\"],\n  [\"16\u00d762=992\", \"88\u00d712=1056\"],\n  [\"73\u00d753=3869\", \"21\u00d751=1071\"],\n  [\"74\u00d714=1036\", \"86\u00d738=3268\"],\n  [\"94\u00d785=7990\", \"97\u00d797=9409\"],\n  [\"35\u00d719=665\", \"48\u00d796=4608\"],\n  [\"45\u00d755=2475\", \"89\u00d764=5696\"],\n  [\"19\u00d733=627\", \"18\u00d777=1386\"],\n  [\"15\u00d755=825\", \"63\u00d789=5607\"],\n  [\"64\u00d775=4800\", \"18\u00d713=234\"],\n  [\"78\u00d735=2730\", \"14\u00d711=154\"],\n  [\"68\u00d784=5712\", \"73\u00d732=2336\"],\n  [\"71\u00d798=6958\", \"73\u00d796=7008\"],\n  [\"53\u00d714=742\", \"61\u00d797=5917\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: update the date line\n# and each two-digit multiplication answer cell in the table.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-08-21 Monday\", \"2023-08-22 Tuesday\"),\n    @(\"41\u00d722=902\", \"75\u00d759=4425\"),\n    @(\"14\u00d750=700\", \"97\u00d732=3104\"),\n    @(\"61\u00d766=4026\", \"75\u00d787=6525\"),\n    @(\"44\u00d726=1144\", \"15\u00d747=705\"),\n    @(\"41\u00d725=1025\", \"20\u00d733=660\"),\n    @(\"87\u00d752=4524\", \"37\u00d762=2294\"),\n    @(\"15\u00d757=855\", \"64\u00d780=5120\"),\n    @(\"28\u00d713=364\", \"53\u00d799=5247\"),\n    @(\"36\u00d728=1008\", \"29\u00d720=580\"),\n    @(\"94\u00d791=8554\", \"95\u00d796=9120\"),\n    @(\"21\u00d789=1869\", \"14\u00d749=686\"),\n    @(\"64\u00d785=5440\", \"92\u00d749=4508\"),\n    @(\"16\u00d762=992\", \"88\u00d712=1056\"),\n    @(\"73\u00d753=3869\", \"21\u00d751=1071\"),\n    @(\"74\u00d714=1036\", \"86\u00d738=3268\"),\n    @(\"94\u00d785=7990\", \"97\u00d797=9409\"),\n    @(\"35\u00d719=665\", \"48\u00d796=4608\"),\n    @(\"45\u00d755=2475\", \"89\u00d764=5696\"),\n    @(\"19\u00d733=627\", \"18\u00d777=1386\"),\n    @(\"15\u00d755=825\", \"63\u00d789=5607\"),\n    @(\"64\u00d775=4800\", \"18\u00d713=234\"),\n    @(\"78\u00d735=2730\", \"14\u00d711=154\"),\n    @(\"68\u00d784=5712\", \"73\u00d732=2336\"),\n    @(\"71\u00d798=6958\", \"73\u00d796=7008\"),\n    @(\"53\u00d714=742\", \"61\u00d797=5917\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
